# Generate Report for Handoff
# - Flip the "Handed back: in sync with en-US" status to "Ready for handoff"
#   on the Overview, zh-cn and de-de sheets.
# - Refresh the associated "Latest HO Xliff Generate Date" / "Latest Handoff
#   Datetime" timestamps.
# - The Status column narrows (was sized for the long "Handed back..."
#   string); re-set it to the width the shorter "Ready for handoff" text
#   needs.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: E2 (zh-cn status), F2 (de-de status), G2 (generate date)
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-27 10:58:58"

# zh-cn sheet: C2 (status), H2 (latest handoff datetime)
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-27 10:58:54"

# de-de sheet: C2 (status), H2 (latest handoff datetime)
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-27 10:58:58"

# Narrow the Status columns to fit the shorter text (was auto-sized for the
# long "Handed back: in sync with en-US" string).
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
